# Add a new weekly price record for "Albahaca" at
# "Feria Lagunitas de Puerto Montt". The new observation (date 2021-12-14,
# serial 44544) is inserted as row 22, pushing the existing rows 22-82
# down to 23-83 (matching the target dimension A1:R83).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 22 (shifts rows 22:82 -> 23:83)
$ws.Rows.Item(22).EntireRow.Insert()

# Populate the newly inserted row 22 with the new data point
$ws.Range("A22").Value2 = 4
$ws.Range("B22").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C22").Value2 = "Los Lagos"
$ws.Range("D22").Value2 = 44544
$ws.Range("E22").Value2 = 10
$ws.Range("F22").Value2 = 100112052
$ws.Range("G22").Value2 = "Albahaca"
$ws.Range("H22").Value2 = "Sin especificar"
$ws.Range("I22").Value2 = "Primera"
$ws.Range("J22").Value2 = 80
$ws.Range("K22").Value2 = 7000
$ws.Range("L22").Value2 = 7000
$ws.Range("M22").Value2 = 7000
$ws.Range("N22").Value2 = "`$/docena de matas"
$ws.Range("O22").Value2 = "Región Metropolitana"
$ws.Range("P22").Value2 = 1167
$ws.Range("Q22").Value2 = 6
$ws.Range("R22").Value2 = "Hortaliza"
